$d = $word.ActiveDocument

# Locate the paragraph that contains the target sentence.
$search = $d.Content
$found = $search.Find.Execute("管理員操作需有限權限控管。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Target sentence not found"
}

$para = $search.Paragraphs(1)
$pRange = $para.Range

# Rebuild the paragraph, splitting the final run into three runs:
#   1) "管理員操作需"            - keeps the original run's rsidRPr, no hint
#   2) "要"                      - new run, rFonts carries hint="eastAsia"
#   3) "有限權限控管。"          - new run, no hint
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="5DC7B16A" w14:textId="77777777" w:rsidR="001E7558" w:rsidRPr="00392DF2" w:rsidRDefault="001E7558" w:rsidP="001E7558" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:pPr><w:pStyle w:val="a9"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:snapToGrid w:val="0"/><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00392DF2"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>資料安全性</w:t></w:r>' +
  '<w:r w:rsidRPr="00392DF2"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>：</w:t></w:r>' +
  '<w:r w:rsidRPr="00392DF2"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>使用者帳號資訊需妥善加密儲存，點數與個人紀錄不得外洩</w:t></w:r>' +
  '<w:r w:rsidRPr="00392DF2"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>，</w:t></w:r>' +
  '<w:r w:rsidRPr="00392DF2"><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>管理員操作需</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>要</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>有限權限控管。</w:t></w:r>' +
  '</w:p>'

$pRange.InsertXML($xml)
